$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# E50 was stored as a text ("inline string") value "590024"; the edit
# converts it to a real number so it matches the rest of the bsecode column.
$ws.Cells.Item(50, 5).Value = 590024

# Append two new rows (51 and 52) of screener data, extending the used
# range from A1:H50 to A1:H52. The bsecode column (E) on these two rows
# keeps its original text formatting, so force the cell to text before
# assigning the numeric-looking string (otherwise it gets auto-converted
# to a number, same as E50 above).
$ws.Cells.Item(51, 1).Value = "25/06/2024 07:44:47"
$ws.Cells.Item(51, 2).Value = 1
$ws.Cells.Item(51, 3).Value = "LODHA"
$ws.Cells.Item(51, 4).Value = "Macrotech Developers Ltd"
$ws.Cells.Item(51, 5).NumberFormat = "@"
$ws.Cells.Item(51, 5).Value = "543287"
$ws.Cells.Item(51, 6).Value = -4.36
$ws.Cells.Item(51, 7).Value = 1481.45
$ws.Cells.Item(51, 8).Value = 351210

$ws.Cells.Item(52, 1).Value = "25/06/2024 07:44:47"
$ws.Cells.Item(52, 2).Value = 2
$ws.Cells.Item(52, 3).Value = "FACT"
$ws.Cells.Item(52, 4).Value = "Fertilizers And Chemicals Travancore Limited"
$ws.Cells.Item(52, 5).NumberFormat = "@"
$ws.Cells.Item(52, 5).Value = "590024"
$ws.Cells.Item(52, 6).Value = -1.18
$ws.Cells.Item(52, 7).Value = 1008.95
$ws.Cells.Item(52, 8).Value = 1599067
